# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the refreshed scrape output (commit: "Update gh-pages
# to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F on each affected sheet.
$updates = @{
    2  = 11833
    3  = 11588
    6  = 1047
    8  = 74
    9  = 48
    11 = 10907
    12 = 4208
    13 = 21
    14 = 15
    16 = 2476
    18 = 60
    20 = 138
    21 = 458
    22 = 11167
    23 = 10973
    29 = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
